$wb = $excel.ActiveWorkbook

$wsReadMe = $wb.Worksheets.Item("ReadMe")
$wsBC = $wb.Worksheets.Item("Biomedical Concepts")

# Update the bc_categories / synonyms descriptions on the ReadMe sheet to
# their pluralized wording.
$wsReadMe.Range("C8").Value = "Biomedical Concept categories for the faciliation of API search and extract"
$wsReadMe.Range("C9").Value = "Biomedical Concept synonyms equivalent to BC short name for the facilitation of API search and extraction"

# Restore the last selection on the ReadMe sheet to C10 (without leaving it
# as the active sheet/tab).
$wsReadMe.Range("C10").Select() | Out-Null
$wsBC.Activate()

# Widen column B ("short_name") slightly on the "Biomedical Concepts" sheet.
$wsBC.Columns.Item(2).ColumnWidth = 42.71
